$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Rename sheet
$ws.Name = "hgdhus"

# Update threshold/statistics values
$ws.Range("J1").Value = 44.12683844566345
$ws.Range("D2").Value = 1862
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 99.89270386266094
$ws.Range("H2").Value = 99.83914209115281
$ws.Range("I2").Value = 0.002680965147453083
$ws.Range("J2").Value = 61.97790503501892
$ws.Range("B3").Value = 2084
$ws.Range("D3").Value = 2083
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 62.81262254714966
$ws.Range("B4").Value = 2590
$ws.Range("D4").Value = 2565
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 6
$ws.Range("G4").Value = 99.76662777129522
$ws.Range("H4").Value = 99.07300115874855
$ws.Range("I4").Value = 0.01166407465007776
$ws.Range("J4").Value = 75.03440403938293
$ws.Range("B5").Value = 2025
$ws.Range("D5").Value = 2021
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 99.75320829220138
$ws.Range("H5").Value = 99.85177865612648
$ws.Range("I5").Value = 0.003946719289590528
$ws.Range("J5").Value = 79.2034797668457
$ws.Range("B6").Value = 1759
$ws.Range("E6").Value = 3
$ws.Range("H6").Value = 99.82935153583618
$ws.Range("I6").Value = 0.005672149744753261
$ws.Range("J6").Value = 63.25737929344177
$ws.Range("B7").Value = 2532
$ws.Range("D7").Value = 2531
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 100
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 102.0119519233704
$ws.Range("B8").Value = 2123
$ws.Range("E8").Value = 0
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 0.0004708097928436911
$ws.Range("J8").Value = 70.40609502792358
$ws.Range("J9").Value = 103.5039761066437
$ws.Range("B10").Value = 1799
$ws.Range("E10").Value = 4
$ws.Range("H10").Value = 99.77753058954394
$ws.Range("I10").Value = 0.00222841225626741
$ws.Range("J10").Value = 85.04934310913086
$ws.Range("B11").Value = 1864
$ws.Range("D11").Value = 1861
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 17
$ws.Range("G11").Value = 99.09478168264111
$ws.Range("H11").Value = 99.89264626945787
$ws.Range("I11").Value = 0.01011176157530601
$ws.Range("J11").Value = 50.19991612434387
$ws.Range("J12").Value = 71.11552119255066
$ws.Range("B13").Value = 2394
$ws.Range("D13").Value = 2392
$ws.Range("F13").Value = 19
$ws.Range("G13").Value = 99.21194525093323
$ws.Range("H13").Value = 99.95821145006268
$ws.Range("I13").Value = 0.008291873963515755
$ws.Range("J13").Value = 91.08571648597717
$ws.Range("B14").Value = 1535
$ws.Range("E14").Value = 0
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 67.71066904067993
$ws.Range("J15").Value = 96.03219509124756
$ws.Range("B16").Value = 1988
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 99.94967287367891
$ws.Range("I16").Value = 0.0005032712632108706
$ws.Range("J16").Value = 109.6916506290436
$ws.Range("D17").Value = 1861
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 99.94629430719657
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 0.0005367686527106817
$ws.Range("J17").Value = 101.3471539020538
$ws.Range("J18").Value = 124.910148859024
$ws.Range("J19").Value = 88.54359149932861
$ws.Range("J20").Value = 65.98450899124146
$ws.Range("B21").Value = 2599
$ws.Range("D21").Value = 2597
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 99.88461538461539
$ws.Range("H21").Value = 99.96150885296382
$ws.Range("J21").Value = 114.2861413955688
$ws.Range("B22").Value = 1940
$ws.Range("D22").Value = 1939
$ws.Range("F22").Value = 23
$ws.Range("G22").Value = 98.82772680937819
$ws.Range("I22").Value = 0.01171676006113092
$ws.Range("J22").Value = 126.755690574646
$ws.Range("B23").Value = 2041
$ws.Range("D23").Value = 2040
$ws.Range("F23").Value = 95
$ws.Range("G23").Value = 95.55035128805621
$ws.Range("I23").Value = 0.04447565543071161
$ws.Range("J23").Value = 81.21774959564209
$ws.Range("B24").Value = 2930
$ws.Range("D24").Value = 2910
$ws.Range("E24").Value = 19
$ws.Range("F24").Value = 69
$ws.Range("G24").Value = 97.68378650553878
$ws.Range("H24").Value = 99.351314441789
$ws.Range("I24").Value = 0.02953020134228188
$ws.Range("J24").Value = 173.7253584861755
$ws.Range("B25").Value = 2641
$ws.Range("D25").Value = 2638
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 99.35969868173258
$ws.Range("H25").Value = 99.92424242424242
$ws.Range("I25").Value = 0.007153614457831325
$ws.Range("J25").Value = 138.5689563751221
$ws.Range("B26").Value = 1855
$ws.Range("D26").Value = 1851
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 99.56966110812265
$ws.Range("H26").Value = 99.83818770226537
$ws.Range("I26").Value = 0.005913978494623656
$ws.Range("J26").Value = 113.9880638122559
$ws.Range("B27").Value = 2945
$ws.Range("D27").Value = 2943
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 99.62762356127286
$ws.Range("H27").Value = 99.96603260869566
$ws.Range("I27").Value = 0.004060913705583757
$ws.Range("J27").Value = 171.7406742572784
$ws.Range("B28").Value = 3005
$ws.Range("E28").Value = 0
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 168.3317763805389
$ws.Range("B29").Value = 2601
$ws.Range("D29").Value = 2599
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 50
$ws.Range("G29").Value = 98.1124952812382
$ws.Range("H29").Value = 99.96153846153847
$ws.Range("I29").Value = 0.01924528301886792
$ws.Range("J29").Value = 146.8044946193695
$ws.Range("J30").Value = 194.5241487026215
$ws.Range("B31").Value = 3250
$ws.Range("D31").Value = 3248
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 99.93846153846154
$ws.Range("H31").Value = 99.96922129886119
$ws.Range("J31").Value = 287.8020353317261
$ws.Range("B32").Value = 2262
$ws.Range("D32").Value = 2248
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 13
$ws.Range("G32").Value = 99.42503317116321
$ws.Range("H32").Value = 99.42503317116321
$ws.Range("I32").Value = 0.01149425287356322
$ws.Range("J32").Value = 176.4919924736023
$ws.Range("B33").Value = 3358
$ws.Range("D33").Value = 3357
$ws.Range("F33").Value = 5
$ws.Range("G33").Value = 99.85127900059489
$ws.Range("I33").Value = 0.001486767766874814
$ws.Range("J33").Value = 228.7169740200043
$ws.Range("D34").Value = 2153
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 100
$ws.Range("H34").Value = 99.9535747446611
$ws.Range("I34").Value = 0.0004642525533890436
$ws.Range("J34").Value = 174.3257346153259
$ws.Range("J35").Value = 242.6634504795074
$ws.Range("B36").Value = 2425
$ws.Range("D36").Value = 2415
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 99.54657873042045
$ws.Range("H36").Value = 99.62871287128714
$ws.Range("I36").Value = 0.008240626287597858
$ws.Range("J36").Value = 160.6088354587555
$ws.Range("B37").Value = 2342
$ws.Range("D37").Value = 2341
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 141
$ws.Range("G37").Value = 94.3190975020145
$ws.Range("H37").Value = 100
$ws.Range("I37").Value = 0.05678614579138139
$ws.Range("J37").Value = 194.821121931076
$ws.Range("B38").Value = 2605
$ws.Range("D38").Value = 2604
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 100
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 135.7912087440491
$ws.Range("B39").Value = 2052
$ws.Range("D39").Value = 2044
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 8
$ws.Range("G39").Value = 99.61013645224172
$ws.Range("H39").Value = 99.65870307167235
$ws.Range("I39").Value = 0.007306380905991232
$ws.Range("J39").Value = 200.9912657737732
$ws.Range("J40").Value = 245.4309468269348
$ws.Range("J41").Value = 145.3222868442535
$ws.Range("B42").Value = 1780
$ws.Range("D42").Value = 1778
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 99.94378864530636
$ws.Range("H42").Value = 99.94378864530636
$ws.Range("I42").Value = 0.001123595505617978
$ws.Range("J42").Value = 117.1968619823456
$ws.Range("D43").Value = 3071
$ws.Range("E43").Value = 6
$ws.Range("F43").Value = 7
$ws.Range("G43").Value = 99.772579597141
$ws.Range("H43").Value = 99.80500487487812
$ws.Range("I43").Value = 0.004222150048717116
$ws.Range("J43").Value = 302.1989524364471
$ws.Range("B44").Value = 2754
$ws.Range("D44").Value = 2752
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 100
$ws.Range("H44").Value = 99.96367598982928
$ws.Range("J44").Value = 253.359103679657
